# Atualização automática de AUGUSTO_PESTANA.xlsx
#
# - Remove the "Desarquivamentos Pendentes" sheet (no longer needed).
# - Normalize the "Paineis DARQ" sheet name to upper case.
# - Normalize the "Recolhimento x Eliminacao" sheet name to upper case
#   (with the accented "ELIMINAÇÃO").

$wb = $excel.ActiveWorkbook

# Avoid the "are you sure you want to delete this sheet" confirmation dialog.
[void]($excel.DisplayAlerts = $false)

# Delete the sheet that is no longer needed.
[void]$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Rename sheets to their new (upper-case) titles.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

[void]($excel.DisplayAlerts = $true)
